$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODOLIST")

# --- Row 3 ---
# O3: "BUG" -> "STABLE", style changes from the red "BUG" look to the
# theme-5 ("En cours"/"READY") look, so copy formatting from a cell that
# already carries that style (G7) before writing the new text.
$ws.Range("G7").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$ws.Range("O3").Value = "STABLE"

# --- Row 5 ---
# M5 / O5 were blank cells that already carry the right formatting
# (style 4 / style 7 respectively), so just fill in the values.
$ws.Range("M5").Value = "Aurélien"
$ws.Range("O5").Value = "Terminé"

# --- Row 6 ---
# J6 / M6 / O6 were blank cells that already carry the right formatting.
$ws.Range("J6").Value = "HallOfFame - level handling"
$ws.Range("M6").Value = "Brice"
$ws.Range("O6").Value = "Terminé"

# --- Row 9 ---
# G9: "Phase d'incubation" -> "En cours", style moves from the red "BUG"
# look to the theme-5 ("En cours") look -> copy formatting from G7.
$ws.Range("G7").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = "En cours"

# --- Row 11 ---
# B11 / E11 were blank cells that already carry the right formatting.
$ws.Range("B11").Value = "Ecran well played (level successful)"
$ws.Range("E11").Value = "Aurélien"
# G11 did not previously hold any value/style -> copy the "Termine" look
# from G3 before setting its text.
$ws.Range("G3").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$ws.Range("G11").Value = "Terminé"

$excel.CutCopyMode = $false

# Match the saved selection/active cell from the diff.
$ws.Range("D9").Select()
